$d = $word.ActiveDocument

# 1. Heading2 "Some Class Name" -> "Robot"
$heading = $d.Paragraphs.Item(1).Range
$heading.Find.Execute("Some Class Name", $true, $false, $false, $false, $false, $true, 1, $false, "Robot", 2)

# 2. Table cell "Class Name" -> "Ro" + bot, with the _GoBack bookmark moved
#    in between the two new runs.
$t = $d.Tables.Item(1)
$cell = $t.Cell(1, 1)
$cellRange = $cell.Range
$cellRange.Text = "Robot"

# Recompute the cell range after the text replacement and split it so the
# bookmark sits right after "Ro".
$cell = $t.Cell(1, 1)
$cellStart = $cell.Range.Start
$bmRange = $d.Range($cellStart + 2, $cellStart + 2)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 3. Merge the split "So" / "me paragraph about this class" runs (which used
#    to carry the _GoBack bookmark) back into a single run of text. Moving
#    the bookmark above already detached it from this location.
$body = $d.Content
$body.Find.ClearFormatting()
$body.Find.Execute("So" + "me paragraph about this class", $true, $false, $false, $false, $false, $true, 1, $false, "Some paragraph about this class", 2)
